# Fruta / hortaliza, semanal
# Insert 3 new weekly-report rows at the top of the Frutilla dataset
# (before the existing row 573), pushing the rest of the table down by
# three rows (573-629 -> 576-632), and populate the 3 new rows with the
# new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before row 573 - this shifts rows 573:629 down
# to 576:632 and extends the used range / dimension automatically.
$ws.Range("A573:A575").EntireRow.Insert()

# Constant columns shared by every record in this dataset block.
$constA = 3
$constB = "Femacal de La Calera"
$constC = "Coquimbo"
$constE = 5
$constF = "Fruta"
$constG = 100101
$constH = "Berries"
$constI = 100112025
$constJ = "Frutilla"
$constK = "Sin especificar"
$constQ = "$/bandeja 7 kilos"
$constR = "Provincia de Melipilla"
$constT = 7

function Fill-Row($r, $fecha, $calidad, $volumen, $pmin, $pmax, $pprom, $pkg) {
    $ws.Cells.Item($r, 1).Value2 = $constA
    $ws.Cells.Item($r, 2).Value2 = $constB
    $ws.Cells.Item($r, 3).Value2 = $constC
    $ws.Cells.Item($r, 4).Value2 = $fecha
    $ws.Cells.Item($r, 5).Value2 = $constE
    $ws.Cells.Item($r, 6).Value2 = $constF
    $ws.Cells.Item($r, 7).Value2 = $constG
    $ws.Cells.Item($r, 8).Value2 = $constH
    $ws.Cells.Item($r, 9).Value2 = $constI
    $ws.Cells.Item($r, 10).Value2 = $constJ
    $ws.Cells.Item($r, 11).Value2 = $constK
    $ws.Cells.Item($r, 12).Value2 = $calidad
    $ws.Cells.Item($r, 13).Value2 = $volumen
    $ws.Cells.Item($r, 14).Value2 = $pmin
    $ws.Cells.Item($r, 15).Value2 = $pmax
    $ws.Cells.Item($r, 16).Value2 = $pprom
    $ws.Cells.Item($r, 17).Value2 = $constQ
    $ws.Cells.Item($r, 18).Value2 = $constR
    $ws.Cells.Item($r, 19).Value2 = $pkg
    $ws.Cells.Item($r, 20).Value2 = $constT
}

Fill-Row 573 45223 "Especial" 80  10000 10000 10000 1429
Fill-Row 574 45223 "Primera"  220 8000  8500  8227  1175
Fill-Row 575 45223 "Segunda"  80  6000  6000  6000  857

# Keep the date column formatted the same way as the rest of the table.
$ws.Range("D573:D575").NumberFormat = $ws.Range("D576").NumberFormat
